$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers (existing): A=ID  B=ServerID  C=Name  D=MaxOnline  E=CpuCount  F=IP  G=Port
# Populate row 2 with the first server's data ("MasterServer_1").
#
# Write order matters for shared-string allocation order (to match the
# target file's sharedStrings.xml index order: 7=000106001, 8=127.0.0.1,
# 9=MasterServer_1), so B2/F2 are written before A2/C2.
$ws.Range("B2").Value = "000106001"
$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("A2").Value = "MasterServer_1"
$ws.Range("C2").Value = "MasterServer_1"
# C2 is a brand-new cell (previously blank/non-existent), so it needs the
# text-format style ("@" = numFmtId 49) explicitly applied to match A2/B2's
# pre-existing style (s="1").
$ws.Range("C2").NumberFormat = "@"
$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 2001

# Move/collapse the selection to H3 (single cell), matching the saved view state.
$ws.Range("H3").Select()
